$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.351.48"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "2.040.02"
$ws.Range("E3").Value = "  +3.96%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.17"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +4.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0807"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.40%  "
$ws.Range("D13").Value = "2.343.45"
$ws.Range("E13").Value = "  +4.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.851"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("D17").Value = "2.037.54"
$ws.Range("E17").Value = "  +4.24%  "
$ws.Range("D18").Value = "37.330.19"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.59%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("E28").Value = "  -4.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.12%  "
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E32").Value = "  +11.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  +12.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.30%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0982"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0217"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.54%  "
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("D47").Value = "1.388.41"
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("E49").Value = "  +14.40%  "
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("D51").Value = "2.235.52"
$ws.Range("E51").Value = "  +4.26%  "
